# i-363-2 | fixed
# Se modificó la forma de importar productos, ahora el almacén se
# selecciona en el popup de importación: la columna "Almacén" (U) se
# quita de la hoja de importación y se reemplaza la fila de ejemplo por
# un nuevo producto de muestra.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nueva fila de producto de ejemplo (reemplaza COCA COLA / TRIPLE COLA).
# B2 se escribe antes que A2 para que el nuevo texto compartido quede en
# el mismo orden que usa el archivo de referencia.
$ws.Range("B2").Value = "BI001"
$ws.Range("A2").Value = "BILLETERA D&G COLOR NEGRO"
$ws.Range("O2").Value = "D&G"
$ws.Range("P2").Value = "BILLETERA D&G COLOR NEGRO"
$ws.Range("Q2").Value = "BILLETERA D&G COLOR NEGRO"
$ws.Range("T2").Value = "B1000001"
$ws.Range("G2").Value = 90
$ws.Range("J2").Value = 70

# La columna "Almacén" (U) ya no se completa en la hoja: ahora el
# almacén se elige en el popup de importación, así que se borra por
# completo (encabezado y dato) en vez de dejarla vacía.
$ws.Range("U1").Value = ""
$ws.Range("U2").Value = ""

# Deja la selección/scroll donde el usuario terminó de revisar el nuevo
# formato de importación.
$ws.Range("B3").Select()
